# Tagsample.pptx edit: turn the escaped-angle-bracket pseudo-tags
# (<gN> ... </gN>) that were typed into the run text back into the
# curly-brace placeholder form ({gN} ... {/gN}) used by the Okapi
# OpenXML encoder tests.
#
# We walk each shape's TextRange one original run (or line break) at a
# time using Characters(start, length) so that run boundaries and their
# run properties (rPr, incl. err="1" spell-check flags) are left
# completely untouched -- only the literal "<" / ">" characters inside
# each run's text are swapped for "{" / "}".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Convert-TagBrackets($TextRange, $Tokens) {
    $pos = 1
    foreach ($tok in $Tokens) {
        if ($tok -eq "BR") {
            # <a:br/> (or paragraph break) - one character, leave as-is.
            $pos = $pos + 1
            continue
        }

        $len = [int]$tok
        if ($len -gt 0) {
            $run = $TextRange.Characters($pos, $len)
            $orig = $run.Text
            $new = $orig.Replace("<", "{").Replace(">", "}")
            if ($new -ne $orig) {
                $run.Text = $new
            }
        }
        $pos = $pos + $len
    }
}

# Shape 1: "Title 1" -- single run "<g0>=andray()</g1>"
$titleShape = $s.Shapes.Item(1)
$titleTokens = @(18)
Convert-TagBrackets $titleShape.TextFrame.TextRange $titleTokens

# Shape 2: "Content Placeholder 2" -- 93 runs spread across 3 soft-wrapped
# lines (separated by <a:br/> elements, each with an empty run before the
# second break of each pair).
$bodyShape = $s.Shapes.Item(2)
$bodyTokens = @(
    23, 16, 16, 16, 24, 17, 33, 18, 18, 18, 26, 17, 33, 18, 18, 18, 26, 17,
    33, 18, 18, 18, 26, 17, 33, 18, 18, 18, 26, 17, 19,
    "BR", 0, "BR",
    25, 18, 18, 18, 26, 17, 33, 18, 18, 18, 26, 17, 33, 18, 18, 18, 26, 18,
    35, 20, 20, 20, 28, 19, 35, 20, 20, 20, 28, 19, 21,
    "BR", 0, "BR",
    27, 20, 20, 20, 28, 19, 35, 20, 20, 20, 28, 19, 35, 20, 20, 20, 28, 19,
    35, 20, 20, 20, 28, 19, 35, 20, 20, 20, 28, 19, 21,
    "BR"
)
Convert-TagBrackets $bodyShape.TextFrame.TextRange $bodyTokens
